# Insert a new data row at row 173 (pushing existing rows 173:274 down to 174:275)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173, shifting rows 173:274 down to 174:275.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(173, 1).Value = 4
$ws.Cells.Item(173, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(173, 3).Value = "Los Lagos"
$ws.Cells.Item(173, 4).Value = 44767
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100108
$ws.Cells.Item(173, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(173, 9).Value = 100108005
$ws.Cells.Item(173, 10).Value = "Piña"
$ws.Cells.Item(173, 11).Value = "Caramelo"
$ws.Cells.Item(173, 12).Value = "Primera"
$ws.Cells.Item(173, 13).Value = 20
$ws.Cells.Item(173, 14).Value = 23000
$ws.Cells.Item(173, 15).Value = 23000
$ws.Cells.Item(173, 16).Value = 23000
$ws.Cells.Item(173, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(173, 18).Value = "Ecuador"
$ws.Cells.Item(173, 19).Value = 1917
$ws.Cells.Item(173, 20).Value = 12
